# Adds a new "2022-Q4" quarter to the 港股/01658-邮储银行 workbook:
#   1. A new fund-holdings detail sheet "2022-Q4" is inserted right after
#      the "总计" (total) sheet, built from a copy of the "2022-Q3" sheet
#      (identical column layout/styling) with its data replaced.
#   2. The "总计" summary sheet gets a new row 2 for 2022-Q4, with all the
#      previously-existing quarters shifted down by one row.

$wb = $excel.ActiveWorkbook
$totalWs = $wb.Worksheets.Item("总计")
$q3Ws = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------
# 1. Build the new "2022-Q4" detail sheet from a copy of "2022-Q3"
#    (same header row / styling), placed right after "总计".
# ---------------------------------------------------------------------
$q3Ws.Copy($null, $totalWs)
$q4Ws = $wb.Worksheets.Item(2)
$q4Ws.Name = "2022-Q4"

# "2022-Q3" has 23 data rows (rows 2-24); "2022-Q4" only needs 21 data
# rows (rows 2-22), so drop the trailing two rows copied along with it.
$q4Ws.Rows.Item(23).Resize(2).Delete()

# Fund code / name / size / position columns are plain-text (so values
# like "008283" keep their leading zero and "90.50" keeps its trailing
# zero instead of being coerced to numbers).
$q4Ws.Range("B2:G22").NumberFormat = "@"

$q4Data = @(
    @(0,"008283","易方达金融行业股票","20.52","87.72","7.28","1.4939",6),
    @(1,"011152","富兰克林国海兴海回报混合","16.09","87.64","4.79","0.7707",6),
    @(2,"008515","富兰克林国海基本面优选混合","12.93","88.93","5.64","0.7293",5),
    @(3,"009983","永赢港股通品质生活慧选混合","9.48","90.50","5.73","0.5432",3),
    @(4,"011468","国富竞争优势三年持有期混合A","13.00","87.96","3.74","0.4862",5),
    @(5,"010365","鹏华港股通中证香港银行投资指数（LOF）C","6.46","94.47","3.54","0.2287",9),
    @(6,"011315","永赢港股通优质成长一年混合","3.53","85.99","5.71","0.2016",3),
    @(7,"011913","华夏永泓一年持有混合A","9.48","32.14","1.14","0.1081",7),
    @(8,"011914","华夏永泓一年持有混合C","9.45","32.14","1.14","0.1077",7),
    @(9,"501025","鹏华港股通中证香港银行投资指数（LOF）A","2.38","94.47","3.54","0.0843",9),
    @(10,"012170","华夏永顺一年持有混合A","8.36","28.10","0.87","0.0727",9),
    @(11,"001703","银华沪港深增长股票A","1.71","93.53","3.41","0.0583",10),
    @(12,"011469","国富竞争优势三年持有期混合C","0.74","87.96","3.74","0.0277",5),
    @(13,"006810","泰康港股通中证香港银行投资指数C","0.58","94.66","3.52","0.0204",9),
    @(14,"006809","泰康港股通中证香港银行投资指数A","0.55","94.66","3.52","0.0194",9),
    @(15,"007751","景顺长城中证沪港深红利成长低波动指数A","0.69","91.50","2.45","0.0169",7),
    @(16,"014364","银华沪港深增长股票C","0.25","93.53","3.41","0.0085",10),
    @(17,"011647","博时港股通红利精选混合A","0.11","77.58","4.32","0.0048",5),
    @(18,"012171","华夏永顺一年持有混合C","0.34","28.10","0.87","0.0030",9),
    @(19,"007760","景顺长城中证沪港深红利成长低波动指数C","0.07","91.50","2.45","0.0017",7),
    @(20,"011648","博时港股通红利精选混合C","0.03","77.58","4.32","0.0013",5)
)

$r = 2
foreach ($row in $q4Data) {
    $q4Ws.Cells.Item($r,1).Value = $row[0]
    $q4Ws.Cells.Item($r,2).Value = $row[1]
    $q4Ws.Cells.Item($r,3).Value = $row[2]
    $q4Ws.Cells.Item($r,4).Value = $row[3]
    $q4Ws.Cells.Item($r,5).Value = $row[4]
    $q4Ws.Cells.Item($r,6).Value = $row[5]
    $q4Ws.Cells.Item($r,7).Value = $row[6]
    $q4Ws.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 2. Insert a 2022-Q4 row into the "总计" summary sheet, shifting the
#    existing quarters (currently rows 2-9) down to rows 3-10.
# ---------------------------------------------------------------------

# Extend the A:D formatting down to the new last row (10) by copying
# the previous last row's formatting before shuffling values around.
$totalWs.Range("A9:D9").Copy($totalWs.Range("A10:D10"))

for ($row = 9; $row -ge 2; $row--) {
    $a = $totalWs.Cells.Item($row, 1).Value()
    $b = $totalWs.Cells.Item($row, 2).Value()
    $c = $totalWs.Cells.Item($row, 3).Value()
    $d = $totalWs.Cells.Item($row, 4).Value()
    $totalWs.Cells.Item($row + 1, 1).Value = $a
    $totalWs.Cells.Item($row + 1, 2).Value = $b
    $totalWs.Cells.Item($row + 1, 3).Value = $c
    $totalWs.Cells.Item($row + 1, 4).Value = $d
}

$totalWs.Cells.Item(2, 1).Value = 0
$totalWs.Cells.Item(2, 2).Value = "2022-Q4"
$totalWs.Cells.Item(2, 3).Value = 21
$totalWs.Cells.Item(2, 4).Value = 4.99
